$d = $word.ActiveDocument

# Update the date line (unique text, safe to use Find/Replace)
$d.Content.Find.Execute("2025-05-11 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-12 Monday", 2) | Out-Null

# Update each answer cell in the table by position (some values repeat,
# so Find/Replace-all would be unsafe -- address cells directly instead)
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "7+31=38"
$t.Cell(1,2).Range.Text = "11+75=86"
$t.Cell(1,3).Range.Text = "43+23=66"
$t.Cell(1,4).Range.Text = "16+54=70"
$t.Cell(1,5).Range.Text = "25+41=66"

$t.Cell(2,1).Range.Text = "44-18=26"
$t.Cell(2,2).Range.Text = "51+24=75"
$t.Cell(2,3).Range.Text = "54+1=55"
$t.Cell(2,4).Range.Text = "39-2=37"
$t.Cell(2,5).Range.Text = "72-14=58"

$t.Cell(3,1).Range.Text = "63-58=5"
$t.Cell(3,2).Range.Text = "55-46=9"
$t.Cell(3,3).Range.Text = "15+81=96"
$t.Cell(3,4).Range.Text = "42+31=73"
$t.Cell(3,5).Range.Text = "79-18=61"

$t.Cell(4,1).Range.Text = "63+12=75"
$t.Cell(4,2).Range.Text = "13+38=51"
$t.Cell(4,3).Range.Text = "23+31=54"
$t.Cell(4,4).Range.Text = "38+38=76"
$t.Cell(4,5).Range.Text = "31+0=31"

$t.Cell(5,1).Range.Text = "30-9=21"
$t.Cell(5,2).Range.Text = "53-7=46"
$t.Cell(5,3).Range.Text = "71-17=54"
$t.Cell(5,4).Range.Text = "15-8=7"
$t.Cell(5,5).Range.Text = "52-9=43"

$t.Cell(6,1).Range.Text = "77-72=5"
$t.Cell(6,2).Range.Text = "1+39=40"
$t.Cell(6,3).Range.Text = "43+6=49"
$t.Cell(6,4).Range.Text = "35-4=31"
$t.Cell(6,5).Range.Text = "60-8=52"

$t.Cell(7,1).Range.Text = "39-36=3"
$t.Cell(7,2).Range.Text = "36-21=15"
$t.Cell(7,3).Range.Text = "18+18=36"
$t.Cell(7,4).Range.Text = "37+51=88"
$t.Cell(7,5).Range.Text = "59-17=42"

$t.Cell(8,1).Range.Text = "69+27=96"
$t.Cell(8,2).Range.Text = "57+15=72"
$t.Cell(8,3).Range.Text = "42-40=2"
$t.Cell(8,4).Range.Text = "10+40=50"
$t.Cell(8,5).Range.Text = "29-29=0"

$t.Cell(9,1).Range.Text = "54-33=21"
$t.Cell(9,2).Range.Text = "89-3=86"
$t.Cell(9,3).Range.Text = "86-43=43"
$t.Cell(9,4).Range.Text = "62-10=52"
$t.Cell(9,5).Range.Text = "89-0=89"

$t.Cell(10,1).Range.Text = "25+32=57"
$t.Cell(10,2).Range.Text = "17+37=54"
$t.Cell(10,3).Range.Text = "54+6=60"
$t.Cell(10,4).Range.Text = "10+76=86"
$t.Cell(10,5).Range.Text = "73+26=99"

$t.Cell(11,1).Range.Text = "67-63=4"
$t.Cell(11,2).Range.Text = "46+25=71"
$t.Cell(11,3).Range.Text = "85-30=55"
$t.Cell(11,4).Range.Text = "0+29=29"
$t.Cell(11,5).Range.Text = "54+32=86"

$t.Cell(12,1).Range.Text = "35-16=19"
$t.Cell(12,2).Range.Text = "19+30=49"
$t.Cell(12,3).Range.Text = "33+9=42"
$t.Cell(12,4).Range.Text = "58-26=32"
$t.Cell(12,5).Range.Text = "50+23=73"

$t.Cell(13,1).Range.Text = "7+7=14"
$t.Cell(13,2).Range.Text = "23+10=33"
$t.Cell(13,3).Range.Text = "96-30=66"
$t.Cell(13,4).Range.Text = "66+10=76"
$t.Cell(13,5).Range.Text = "51+10=61"

$t.Cell(14,1).Range.Text = "50-31=19"
$t.Cell(14,2).Range.Text = "6+13=19"
$t.Cell(14,3).Range.Text = "77-38=39"
$t.Cell(14,4).Range.Text = "20-10=10"
$t.Cell(14,5).Range.Text = "49-22=27"

$t.Cell(15,1).Range.Text = "0+0=0"
$t.Cell(15,2).Range.Text = "57-14=43"
$t.Cell(15,3).Range.Text = "3+71=74"
$t.Cell(15,4).Range.Text = "16+29=45"
$t.Cell(15,5).Range.Text = "78-31=47"

$t.Cell(16,1).Range.Text = "63-9=54"
$t.Cell(16,2).Range.Text = "69-14=55"
$t.Cell(16,3).Range.Text = "35-31=4"
$t.Cell(16,4).Range.Text = "80-34=46"
$t.Cell(16,5).Range.Text = "67-47=20"

$t.Cell(17,1).Range.Text = "96-59=37"
$t.Cell(17,2).Range.Text = "87-46=41"
$t.Cell(17,3).Range.Text = "92-70=22"
$t.Cell(17,4).Range.Text = "9+29=38"
$t.Cell(17,5).Range.Text = "42+42=84"

$t.Cell(18,1).Range.Text = "24+63=87"
$t.Cell(18,2).Range.Text = "10+31=41"
$t.Cell(18,3).Range.Text = "34-9=25"
$t.Cell(18,4).Range.Text = "80-39=41"
$t.Cell(18,5).Range.Text = "44+45=89"

$t.Cell(19,1).Range.Text = "47-21=26"
$t.Cell(19,2).Range.Text = "4+4=8"
$t.Cell(19,3).Range.Text = "20+78=98"
$t.Cell(19,4).Range.Text = "54-34=20"
$t.Cell(19,5).Range.Text = "96-69=27"

$t.Cell(20,1).Range.Text = "62-14=48"
$t.Cell(20,2).Range.Text = "43+2=45"
$t.Cell(20,3).Range.Text = "74-73=1"
$t.Cell(20,4).Range.Text = "51+31=82"
$t.Cell(20,5).Range.Text = "46-40=6"
